$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.038164951147029
$ws.Range("D2").Value = 1.040809318403496
$ws.Range("E2").Value = 1.044820413392423
$ws.Range("F2").Value = 1.050930401740804
$ws.Range("I2").Value = 1.026109458634859
$ws.Range("J2").Value = 1.043263872156114
$ws.Range("K2").Value = 1.043590614427709
$ws.Range("L2").Value = 1.047590388517168
$ws.Range("M2").Value = 1.053683312125302
$ws.Range("N2").Value = 1.018228144688219
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.040627658789585
$ws.Range("D3").Value = 1.043186676703034
$ws.Range("E3").Value = 1.047083830186833
$ws.Range("F3").Value = 1.053488880866137
$ws.Range("I3").Value = 1.02637992601388
$ws.Range("J3").Value = 1.045363680007396
$ws.Range("K3").Value = 1.045774358955608
$ws.Range("L3").Value = 1.049661339990159
$ws.Range("M3").Value = 1.056049848430944
$ws.Range("N3").Value = 1.018927408508906
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.042203379657983
$ws.Range("D4").Value = 1.044707622238164
$ws.Range("E4").Value = 1.048529366801978
$ws.Range("F4").Value = 1.055116559133274
$ws.Range("I4").Value = 1.026543888854601
$ws.Range("J4").Value = 1.046704720390504
$ws.Range("K4").Value = 1.047169763879089
$ws.Range("L4").Value = 1.050982053555573
$ws.Range("M4").Value = 1.057553122242829
$ws.Range("N4").Value = 1.019373846371572
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.042861638184914
$ws.Range("D5").Value = 1.045342957738133
$ws.Range("E5").Value = 1.049132598382612
$ws.Range("F5").Value = 1.05579428485791
$ws.Range("I5").Value = 1.026610191799708
$ws.Range("J5").Value = 1.04726434375271
$ws.Range("K5").Value = 1.047752256254797
$ws.Range("L5").Value = 1.051532737615271
$ws.Range("M5").Value = 1.058178494436439
$ws.Range("N5").Value = 1.019560112016507
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.042971920458508
$ws.Range("D6").Value = 1.045449397162436
$ws.Range("E6").Value = 1.049233623853265
$ws.Range("F6").Value = 1.055907697043013
$ws.Range("I6").Value = 1.02662117079665
$ws.Range("J6").Value = 1.047358065927589
$ws.Range("K6").Value = 1.047849819117854
$ws.Range("L6").Value = 1.051624935804155
$ws.Range("M6").Value = 1.058283112931762
$ws.Range("N6").Value = 1.019591304529406
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.042212191625282
$ws.Range("D7").Value = 1.044716127500807
$ws.Range("E7").Value = 1.048537444669604
$ws.Range("F7").Value = 1.055125640529438
$ws.Range("I7").Value = 1.026544785097209
$ws.Range("J7").Value = 1.046712214305433
$ws.Range("K7").Value = 1.047177563312967
$ws.Range("L7").Value = 1.050989429568031
$ws.Range("M7").Value = 1.057561504306787
$ws.Range("N7").Value = 1.019376340793643
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.039000979271745
$ws.Range("D8").Value = 1.041616406692889
$ws.Range("E8").Value = 1.045589343048241
$ws.Range("F8").Value = 1.051800883005747
$ws.Range("I8").Value = 1.026203162998321
$ws.Range("J8").Value = 1.043977222142336
$ws.Range("K8").Value = 1.044332321518875
$ws.Range("L8").Value = 1.048294330428125
$ws.Range("M8").Value = 1.054488964790887
$ws.Range("N8").Value = 1.018465730105606
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.033201716605394
$ws.Range("D9").Value = 1.036017226497309
$ws.Range("E9").Value = 1.040244428151594
$ws.Range("F9").Value = 1.045723803230082
$ws.Range("I9").Value = 1.025515695887886
$ws.Range("J9").Value = 1.039018604891574
$ws.Range("K9").Value = 1.039179748080983
$ws.Range("L9").Value = 1.043393221942395
$ws.Range("M9").Value = 1.048854983608471
$ws.Range("N9").Value = 1.016813628746163
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.029234732089296
$ws.Range("D10").Value = 1.032186297733158
$ws.Range("E10").Value = 1.036574271257841
$ws.Range("F10").Value = 1.041517753722413
$ws.Range("I10").Value = 1.024998575290955
$ws.Range("J10").Value = 1.035613639160288
$ws.Range("K10").Value = 1.035645577783842
$ws.Range("L10").Value = 1.040017819282708
$ws.Range("M10").Value = 1.044943753162269
$ws.Range("N10").Value = 1.015678415063755
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.02749160917659
$ws.Range("D11").Value = 1.030502774839629
$ws.Range("E11").Value = 1.034958273568809
$ws.Range("F11").Value = 1.039657949519038
$ws.Range("I11").Value = 1.024760377813686
$ws.Range("J11").Value = 1.034114388838487
$ws.Range("K11").Value = 1.034090377877333
$ws.Range("L11").Value = 1.038529243278389
$ws.Range("M11").Value = 1.04321154520811
$ws.Range("N11").Value = 1.015178390343789
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.026840189214692
$ws.Range("D12").Value = 1.029873600928168
$ws.Range("E12").Value = 1.034353865616074
$ws.Range("F12").Value = 1.038961174984526
$ws.Range("I12").Value = 1.024669724432423
$ws.Range("J12").Value = 1.033553643278792
$ws.Range("K12").Value = 1.033508847369891
$ws.Range("L12").Value = 1.037972138843503
$ws.Range("M12").Value = 1.04256216273913
$ws.Range("N12").Value = 1.014991346401676
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.026980101776533
$ws.Range("D13").Value = 1.03000873667378
$ws.Range("E13").Value = 1.034483703308304
$ws.Range("F13").Value = 1.039110907997642
$ws.Range("I13").Value = 1.024689268885132
$ws.Range("J13").Value = 1.033674101638095
$ws.Range("K13").Value = 1.033633764334813
$ws.Range("L13").Value = 1.038091830893908
$ws.Range("M13").Value = 1.042701730015849
$ws.Range("N13").Value = 1.015031528019291
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.02743784389065
$ws.Range("D14").Value = 1.030450846188991
$ws.Range("E14").Value = 1.03490839853643
$ws.Range("F14").Value = 1.039600476511609
$ws.Range("I14").Value = 1.024752928956501
$ws.Range("J14").Value = 1.034068116835472
$ws.Range("K14").Value = 1.034042387854416
$ws.Range("L14").Value = 1.038483278883165
$ws.Range("M14").Value = 1.043157989663534
$ws.Range("N14").Value = 1.015162956246298
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.027719346923633
$ws.Range("D15").Value = 1.030722731910796
$ws.Range("E15").Value = 1.035169512818985
$ws.Range("F15").Value = 1.039901320851213
$ws.Range("I15").Value = 1.024791862702125
$ws.Range("J15").Value = 1.034310367803809
$ws.Range("K15").Value = 1.034293639053373
$ws.Range("L15").Value = 1.038723905070373
$ws.Range("M15").Value = 1.043438311130418
$ws.Range("N15").Value = 1.015243758366163
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.029349871335767
$ws.Range("D16").Value = 1.032297496397267
$ws.Range("E16").Value = 1.036680944304857
$ws.Range("F16").Value = 1.041640356330642
$ws.Range("I16").Value = 1.025014080218598
$ws.Range("J16").Value = 1.035712605240249
$ws.Range("K16").Value = 1.035748256791886
$ws.Range("L16").Value = 1.040116031778787
$ws.Range("M16").Value = 1.045057886864113
$ws.Range("N16").Value = 1.015711418237912
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.030365766979445
$ws.Range("D17").Value = 1.033278602530911
$ws.Range("E17").Value = 1.037621762488688
$ws.Range("F17").Value = 1.042720768816784
$ws.Range("I17").Value = 1.025149626980564
$ws.Range("J17").Value = 1.036585448603547
$ws.Range("K17").Value = 1.036653954486094
$ws.Range("L17").Value = 1.040981960658647
$ws.Range("M17").Value = 1.046063353236726
$ws.Range("N17").Value = 1.016002473680908
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.030955880926487
$ws.Range("D18").Value = 1.033848490300784
$ws.Range("E18").Value = 1.038167949459907
$ws.Range("F18").Value = 1.043347246724527
$ws.Range("I18").Value = 1.025227313404832
$ws.Range("J18").Value = 1.03709217141581
$ws.Range("K18").Value = 1.037179840529026
$ws.Range("L18").Value = 1.041484445483011
$ws.Range("M18").Value = 1.046646109878099
$ws.Range("N18").Value = 1.016171426995517
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.031156684425068
$ws.Range("D19").Value = 1.034042408241668
$ws.Range("E19").Value = 1.038353751913196
$ws.Range("F19").Value = 1.043560235440894
$ws.Range("I19").Value = 1.025253570021934
$ws.Range("J19").Value = 1.037264548621066
$ws.Range("K19").Value = 1.037358751944913
$ws.Range("L19").Value = 1.041655343016001
$ws.Range("M19").Value = 1.046844189481906
$ws.Range("N19").Value = 1.016228898756906
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.030257024385875
$ws.Range("D20").Value = 1.033173585666706
$ws.Range("E20").Value = 1.03752108899719
$ws.Range("F20").Value = 1.04260523554627
$ws.Range("I20").Value = 1.025135226607334
$ws.Range("J20").Value = 1.036492049029094
$ws.Range("K20").Value = 1.036557029984263
$ws.Range("L20").Value = 1.0408893241651
$ws.Range("M20").Value = 1.045955861597446
$ws.Range("N20").Value = 1.015971330707055
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.027303160276372
$ws.Range("D21").Value = 1.030320762966665
$ws.Range("E21").Value = 1.034783452124751
$ws.Range("F21").Value = 1.03945647671229
$ws.Range("I21").Value = 1.024734242983458
$ws.Range("J21").Value = 1.033952196647944
$ws.Range("K21").Value = 1.033922165975773
$ws.Range("L21").Value = 1.038368123648726
$ws.Range("M21").Value = 1.043023798520807
$ws.Range("N21").Value = 1.015124290471257
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.025423045780144
$ws.Range("D22").Value = 1.028504804245656
$ws.Range("E22").Value = 1.033038086425618
$ws.Range("F22").Value = 1.037442152470991
$ws.Range("I22").Value = 1.024469520953183
$ws.Range("J22").Value = 1.032332910412202
$ws.Range("K22").Value = 1.03224312486959
$ws.Range("L22").Value = 1.036758688852479
$ws.Range("M22").Value = 1.041145706328527
$ws.Range("N22").Value = 1.014584107738292
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.026421946528161
$ws.Range("D23").Value = 1.029469633742476
$ws.Range("E23").Value = 1.0339656672844
$ws.Range("F23").Value = 1.038513319371676
$ws.Range("I23").Value = 1.024611060964565
$ws.Range("J23").Value = 1.033193487533395
$ws.Range("K23").Value = 1.033135381608058
$ws.Range("L23").Value = 1.03761422285513
$ws.Range("M23").Value = 1.042144653059414
$ws.Range("N23").Value = 1.014871204495252
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.03030616801726
$ws.Range("D24").Value = 1.033221045597856
$ws.Range("E24").Value = 1.03756658696772
$ws.Range("F24").Value = 1.042657451506902
$ws.Range("I24").Value = 1.025141737765313
$ws.Range("J24").Value = 1.036534259660351
$ws.Range("K24").Value = 1.036600833378792
$ws.Range("L24").Value = 1.040931190625763
$ws.Range("M24").Value = 1.046004443911444
$ws.Range("N24").Value = 1.015985405390486
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.034718259701583
$ws.Range("D25").Value = 1.037481590897951
$ws.Range("E25").Value = 1.041644593787391
$ws.Range("F25").Value = 1.047321536951384
$ws.Range("I25").Value = 1.025703683860893
$ws.Range("J25").Value = 1.040317585393604
$ws.Range("K25").Value = 1.040528844780276
$ws.Range("L25").Value = 1.044678867751517
$ws.Range("M25").Value = 1.05033829216949
$ws.Range("N25").Value = 1.017246552536839
